$wb = $excel.ActiveWorkbook

# --- ALC row 17 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H17").Value = 1953.7142
$ws.Range("I17").Value = 500
$ws.Range("J17").Value = 2065.5386
$ws.Range("K17").Value = 1500
$ws.Range("L17").Value = 6196.6158
$ws.Range("M17").Value = -1332
$ws.Range("N17").Value = -6532.6158

# --- ALC row 86 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H86").Value = 10000
$ws.Range("J86").Value = 10000
$ws.Range("L86").Value = 10000
$ws.Range("N86").Value = -12246

# --- ALC row 89 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H89").Value = 10000
$ws.Range("J89").Value = 10000
$ws.Range("L89").Value = 50000
$ws.Range("N89").Value = -61232

# --- ALC row 112 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H112").Value = 1960.9286
$ws.Range("I112").Value = 1863.3334
$ws.Range("J112").Value = 1987.5454
$ws.Range("K112").Value = 5590.0002
$ws.Range("L112").Value = 5962.6362
$ws.Range("M112").Value = -4482.0002
$ws.Range("N112").Value = -8178.6362

# --- ALC row 137 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H137").Value = 2555.0625
$ws.Range("I137").Value = 1730.3334
$ws.Range("J137").Value = 3049.9
$ws.Range("K137").Value = 5191.0002
$ws.Range("L137").Value = 9149.700000000001
$ws.Range("M137").Value = -2641.0002
$ws.Range("N137").Value = -14249.7

# --- ALC row 138 ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H138").Value = 2965.4546
$ws.Range("J138").Value = 4586.8
$ws.Range("L138").Value = 13760.4
$ws.Range("N138").Value = -24040.4

# --- ARM row 2 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H2").Value = 1582.9333
$ws.Range("I2").Value = 980.7857
$ws.Range("K2").Value = 980.7857
$ws.Range("M2").Value = -867.7857

# --- ARM row 28 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H28").Value = 13647
$ws.Range("J28").Value = 4000
$ws.Range("L28").Value = 4000
$ws.Range("N28").Value = -4384

# --- ARM row 61 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H61").Value = 5511.933
$ws.Range("I61").Value = 3297.6667
$ws.Range("K61").Value = 3297.6667
$ws.Range("M61").Value = -3085.6667

# --- ARM row 99 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H99").Value = 13647
$ws.Range("J99").Value = 4000
$ws.Range("L99").Value = 4000
$ws.Range("N99").Value = -9990

# --- ARM row 102 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H102").Value = 3126.4736
$ws.Range("I102").Value = 1643.0714
$ws.Range("K102").Value = 1643.0714
$ws.Range("M102").Value = -21.07140000000004

# --- ARM row 116 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H116").Value = 1582.9333
$ws.Range("I116").Value = 980.7857
$ws.Range("K116").Value = 980.7857
$ws.Range("M116").Value = 1313.2143

# --- ARM row 124 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H124").Value = 25000
$ws.Range("J124").Value = 25000
$ws.Range("L124").Value = 25000
$ws.Range("N124").Value = -34820

# --- ARM row 132 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H132").Value = 1717.8125
$ws.Range("I132").Value = 1784.9286
$ws.Range("J132").Value = 1248
$ws.Range("K132").Value = 5354.7858
$ws.Range("L132").Value = 3744
$ws.Range("M132").Value = -2824.7858
$ws.Range("N132").Value = -8804

# --- ARM row 136 ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H136").Value = 5511.933
$ws.Range("I136").Value = 3297.6667
$ws.Range("K136").Value = 9893.000100000001
$ws.Range("M136").Value = -7343.000100000001

# --- BSM row 3 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H3").Value = 1582.9333
$ws.Range("I3").Value = 980.7857
$ws.Range("K3").Value = 980.7857
$ws.Range("M3").Value = -866.7857

# --- BSM row 11 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H11").Value = 693.6
$ws.Range("I11").Value = 117
$ws.Range("J11").Value = 3000
$ws.Range("K11").Value = 117
$ws.Range("L11").Value = 3000
$ws.Range("M11").Value = 23
$ws.Range("N11").Value = -3280

# --- BSM row 26 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H26").Value = 15709.857
$ws.Range("I26").Value = 16828.166
$ws.Range("J26").Value = 9000
$ws.Range("K26").Value = 16828.166
$ws.Range("L26").Value = 9000
$ws.Range("M26").Value = -16536.166
$ws.Range("N26").Value = -9584

# --- BSM row 96 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H96").Value = 14600.4
$ws.Range("I96").Value = 15111.556
$ws.Range("K96").Value = 15111.556
$ws.Range("M96").Value = -12365.556

# --- BSM row 99 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H99").Value = 1969.375
$ws.Range("I99").Value = 1679.4286
$ws.Range("K99").Value = 1679.4286
$ws.Range("M99").Value = -181.4286

# --- BSM row 105 ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H105").Value = 1554.1428
$ws.Range("I105").Value = 1489.7
$ws.Range("K105").Value = 1489.7
$ws.Range("M105").Value = 257.3

# --- CRP row 74 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H74").Value = 57524.168
$ws.Range("J74").Value = 57524.168
$ws.Range("L74").Value = 57524.168
$ws.Range("N74").Value = -59272.168

# --- CRP row 77 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H77").Value = 57524.168
$ws.Range("J77").Value = 57524.168
$ws.Range("L77").Value = 172572.504
$ws.Range("N77").Value = -181308.504

# --- CRP row 99 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H99").Value = 1883.3334
$ws.Range("I99").Value = 1883.3334
$ws.Range("K99").Value = 1883.3334
$ws.Range("M99").Value = -385.3334

# --- CRP row 126 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H126").Value = 1883.3334
$ws.Range("I126").Value = 1883.3334
$ws.Range("K126").Value = 5650.0002
$ws.Range("M126").Value = -3180.0002

# --- CRP row 132 ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H132").Value = 2972.25
$ws.Range("I132").Value = 2972.25
$ws.Range("K132").Value = 8916.75
$ws.Range("M132").Value = -6386.75

# --- GSM row 95 ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H95").Value = 25299
$ws.Range("J95").Value = 25299
$ws.Range("L95").Value = 25299
$ws.Range("N95").Value = -30791

# --- LTW row 22 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 956
$ws.Range("I22").Value = 1144.5
$ws.Range("J22").Value = 830.3333
$ws.Range("K22").Value = 1144.5
$ws.Range("L22").Value = 830.3333
$ws.Range("M22").Value = -849.5
$ws.Range("N22").Value = -1420.3333

# --- LTW row 27 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H27").Value = 956
$ws.Range("I27").Value = 1144.5
$ws.Range("J27").Value = 830.3333
$ws.Range("K27").Value = 1144.5
$ws.Range("L27").Value = 830.3333
$ws.Range("M27").Value = -1037.5
$ws.Range("N27").Value = -1044.3333

# --- LTW row 46 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H46").Value = 1000
$ws.Range("I46").Value = 0
$ws.Range("J46").Value = 1000
$ws.Range("K46").Value = 0
$ws.Range("L46").Value = 1000
$ws.Range("M46").ClearContents()
$ws.Range("N46").Value = -1376

# --- LTW row 55 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H55").Value = 1058.625
$ws.Range("I55").Value = 946.8182
$ws.Range("J55").Value = 1304.6
$ws.Range("K55").Value = 946.8182
$ws.Range("L55").Value = 1304.6
$ws.Range("M55").Value = -773.8182
$ws.Range("N55").Value = -1650.6

# --- LTW row 100 ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H100").Value = 5059.778
$ws.Range("I100").Value = 2589.6667
$ws.Range("K100").Value = 2589.6667
$ws.Range("M100").Value = -2048.6667

# --- WVR row 12 ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H12").Value = 4329
$ws.Range("J12").Value = 3993.5
$ws.Range("L12").Value = 3993.5
$ws.Range("N12").Value = -4277.5
